$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two row re-orderings)
# captured from the latest GitHub Actions scrape.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.923.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.256.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.95%  '
$ws.Range("E7").Value = '  -0.92%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.511'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0784'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.27%  '
$ws.Range("E13").Value = '  -1.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.601.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.252.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '46.920.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.44%  '
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0935'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.95%  '
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E26").Value = '  -1.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '42.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '145.91'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.33%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.21'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.29%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0771'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.83%  '
$ws.Range("E36").Value = '  +11.53%  '
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '16.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +20.05%  '
$ws.Range("E39").Value = '  -3.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.15%  '
$ws.Range("E41").Value = '  -4.54%  '
$ws.Range("E42").Value = '  -1.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.93'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +20.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.757.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("B47").Value = 'ordi'
$ws.Range("C47").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '71.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.42%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.186'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '93.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.92%  '
